$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-BottomThin($r) {
  $b = $r.Borders.Item(9)
  $b.ColorIndex = 1
  $b.Weight = 2
  $b.LineStyle = 1
}

function Set-TopBottomThin($r) {
  $b9 = $r.Borders.Item(9)
  $b9.ColorIndex = 1
  $b9.Weight = 2
  $b9.LineStyle = 1
  $b8 = $r.Borders.Item(8)
  $b8.ColorIndex = 1
  $b8.Weight = 2
  $b8.LineStyle = 1
}

# ---- Row 22: restyle existing cells + add empty A22 ----
$a22 = $ws.Range("A22")
$a22.WrapText = $true
Set-BottomThin $a22

$b22 = $ws.Range("B22")
$b22.WrapText = $true
Set-BottomThin $b22

foreach ($addr in @("C22","D22","E22")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
  $r.Font.Size = 8
  Set-BottomThin $r
}

# ---- Prepare formatting for the three new rows first ----
$ws.Rows.Item(23).RowHeight = 43.2
$ws.Rows.Item(24).RowHeight = 43.2
$ws.Rows.Item(25).RowHeight = 43.2

foreach ($addr in @("A23","B23")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
  Set-TopBottomThin $r
}
foreach ($addr in @("C23","D23","E23")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
  $r.Font.Size = 8
  Set-TopBottomThin $r
}

foreach ($addr in @("A24","B24")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
  Set-TopBottomThin $r
}
foreach ($addr in @("C24","D24","E24")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
  $r.Font.Size = 8
  Set-TopBottomThin $r
}

$a25 = $ws.Range("A25")
$a25.WrapText = $true

foreach ($addr in @("B25")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
  Set-TopBottomThin $r
}
foreach ($addr in @("C25","D25","E25")) {
  $r = $ws.Range($addr)
  $r.WrapText = $true
  $r.Font.Size = 8
  Set-TopBottomThin $r
}

# ---- Now write values in the exact order the shared-string table was built ----
# Row 23: English text first, then filename, then Russian, then transliteration
$ws.Range("C23").Value = " HEY, [CS:N]Diglett[CR]![K] We\'ve been away\nfrom this! We\'d better FOCUS!"
$ws.Range("A23").Value = "SCRIPT/G01P04A/um1101.ssb"
$ws.Range("D23").Value = " ЭЙ, [CS:N]Диглетт[CR]![K] Что-то мы\nотвлеклись! Нам нужно СОСРЕДОТОЧИТЬСЯ!"
$ws.Range("E23").Value = " ÜÊ, [CS:N]Äéãìåóó[CR]![K] Œóï-óï íú\nïóâìåëìéòû! Îàí îôçîï ÒÏÒÑÅÄÏÓÏŒÉÓÛÒŸ!"
$ws.Range("B23").Value = 612

# Row 24: filename, then Russian, then transliteration (English reuses row 9's string)
$ws.Range("A24").Value = "SCRIPT/G01P04A/um1103.ssb"
$ws.Range("C24").Value = " HEY, [CS:N]Diglett[CR]! It\'s a new day!\nLet\'s buckle down and FOCUS!"
$ws.Range("D24").Value = " ЭЙ, [CS:N]Диглетт[CR]! Наступил новый\nдень! Пора браться за дело и\nСОСРЕДОТОЧИТЬСЯ!"
$ws.Range("E24").Value = " ÜÊ, [CS:N]Äéãìåóó[CR]! Îàòóôðéì îïâúê\näåîû! Ðïñà áñàóûòÿ èà äåìï é\nÒÏÒÑÅÄÏÓÏŒÉÓÛÒŸ!"
$ws.Range("B24").Value = 593

# Row 25: filename only is new; B/C/D/E reuse existing shared strings
$ws.Range("A25").Value = "SCRIPT/G01P04A/um1105.ssb"
$ws.Range("C25").Value = " HEY, [CS:N]Diglett[CR]! It\'s a new day!\nLet\'s buckle down and FOCUS!"
$ws.Range("D25").Value = " ЭЙ, [CS:N]Диглетт[CR]! Наступил новый\nдень! Пора браться за дело и\nСОСРЕДОТОЧИТЬСЯ!"
$ws.Range("E25").Value = " ÜÊ, [CS:N]Äéãìåóó[CR]! Îàòóôðéì îïâúê\näåîû! Ðïñà áñàóûòÿ èà äåìï é\nÒÏÒÑÅÄÏÓÏŒÉÓÛÒŸ!"
$ws.Range("B25").Value = 593

# ---- Selection / view ----
$ws.Range("C24").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "edit complete"
